$d = $word.ActiveDocument

# Locate the end of the existing last paragraph's text ("...function properly.")
$splitPoint = $d.Content
$splitPoint.Find.Execute("function properly.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint.Collapse(0)

# Split the paragraph here: a new paragraph is created after this point, and
# anything that was anchored at this position (the _GoBack bookmark) travels
# into the new paragraph rather than staying with the old one.
$splitPoint.InsertBefore([char]13)

# The new paragraph is the one right after the paragraph we just split -
# i.e. now the very last paragraph in the document.
$newPara = $d.Paragraphs.Last

# Build the new log entry text before the (now relocated) bookmark.
$insertPoint = $newPara.Range
$insertPoint.Collapse(1)

$insertPoint.InsertBefore("20/07/2017 19:39")
$insertPoint.Collapse(0)

$insertPoint.InsertBefore("`tCompleted Advanced task ‘")
$insertPoint.Collapse(0)

$insertPoint.InsertBefore("Simplify JavaScript")
$insertPoint.Collapse(0)

$insertPoint.InsertBefore("’. Converted populate functions into ")
$insertPoint.Collapse(0)

$insertPoint.InsertBefore("one")
$insertPoint.Collapse(0)

$insertPoint.InsertBefore(" function ")
$insertPoint.Collapse(0)

$insertPoint.InsertBefore("that takes a string to identify each element.")
$insertPoint.Collapse(0)

# Finally, append two trailing spaces after the bookmark, at the very end of
# the (new, now-last) paragraph.
$tail = $newPara.Range
$tail.Collapse(0)
$tail.InsertAfter("  ")
